$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the header text casing: "retarder" -> "Retarder"
$ws.Range("A1").Value = "Retarder"

# 2. Header (A1:D1) now uses the Times New Roman 12pt font used by the data grid
$header = $ws.Range("A1:D1")
$header.Font.Name = "Times New Roman"
$header.Font.Size = 12
$ws.Rows.Item(1).RowHeight = 15

# 3. Column B is widened to fit the new "extra ordinary ray" labels
$ws.Columns.Item(2).ColumnWidth = 16

# 4. Add 13 new rows (6-18) below the matrix, columns A and B, formatted
#    like the rest of the sheet (centered, Times New Roman 12pt) but with
#    no border, for the extra ordinary ray polarimetric acquisitions.
$extra = $ws.Range("A6:B18")
$extra.Font.Name = "Times New Roman"
$extra.Font.Size = 12
$extra.HorizontalAlignment = -4108
$extra.VerticalAlignment = -4107
$extra.Borders.LineStyle = -4142
$ws.Range("A6:A18").RowHeight = 15

# 5. Move the active selection the way the author left it
$ws.Range("C11").Select()
